# Auto-generated edit script: refresh Ridill Profits market-data columns (H-N)
# per sheet, matching the scheduled-runner data refresh described in the commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1138.3572
$ws.Range("I92").Value = 1149.0769
$ws.Range("J92").Value = 999
$ws.Range("K92").Value = 1149.0769
$ws.Range("L92").Value = 999
$ws.Range("M92").Value = 98.92309999999998
$ws.Range("N92").Value = -3495
$ws.Range("H95").Value = 39624
$ws.Range("J95").Value = 39624
$ws.Range("L95").Value = 39624
$ws.Range("N95").Value = -45116
$ws.Range("H96").Value = 506
$ws.Range("I96").Value = 364.3
$ws.Range("J96").Value = 1214.5
$ws.Range("K96").Value = 1092.9
$ws.Range("L96").Value = 3643.5
$ws.Range("M96").Value = 280.0999999999999
$ws.Range("N96").Value = -6389.5
$ws.Range("H97").Value = 2387.5
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 2387.5
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 7162.5
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -8154.5
$ws.Range("H103").Value = 62502988
$ws.Range("I103").Value = 333334340
$ws.Range("J103").Value = 3446.1538
$ws.Range("K103").Value = 1000003020
$ws.Range("L103").Value = 10338.4614
$ws.Range("M103").Value = -1000002434
$ws.Range("N103").Value = -11510.4614
$ws.Range("H121").Value = 1666.6666
$ws.Range("J121").Value = 1666.6666
$ws.Range("L121").Value = 4999.9998
$ws.Range("N121").Value = -8493.9998
$ws.Range("H141").Value = 1558.3636
$ws.Range("I141").Value = 1319.6842
$ws.Range("J141").Value = 3070
$ws.Range("K141").Value = 3959.0526
$ws.Range("L141").Value = 9210
$ws.Range("M141").Value = 1220.9474
$ws.Range("N141").Value = -19570

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8019552
$ws.Range("I32").Value = 1900403.2
$ws.Range("J32").Value = 41674870
$ws.Range("K32").Value = 1900403.2
$ws.Range("L32").Value = 41674870
$ws.Range("M32").Value = -1900116.2
$ws.Range("N32").Value = -41675444
$ws.Range("H94").Value = 33000
$ws.Range("J94").Value = 33000
$ws.Range("L94").Value = 33000
$ws.Range("N94").Value = -34802
$ws.Range("H97").Value = 483.45834
$ws.Range("I97").Value = 478.94446
$ws.Range("J97").Value = 497
$ws.Range("K97").Value = 478.94446
$ws.Range("L97").Value = 497
$ws.Range("M97").Value = 17.05554000000001
$ws.Range("N97").Value = -1489
$ws.Range("H102").Value = 2889.9092
$ws.Range("I102").Value = 3187.375
$ws.Range("J102").Value = 2096.6667
$ws.Range("K102").Value = 3187.375
$ws.Range("L102").Value = 2096.6667
$ws.Range("M102").Value = -1565.375
$ws.Range("N102").Value = -5340.6667
$ws.Range("H132").Value = 9488755
$ws.Range("I132").Value = 10105021
$ws.Range("J132").Value = 6946656
$ws.Range("K132").Value = 30315063
$ws.Range("L132").Value = 20839968
$ws.Range("M132").Value = -30312533
$ws.Range("N132").Value = -20845028

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1157.2222
$ws.Range("I99").Value = 892.7143
$ws.Range("K99").Value = 892.7143
$ws.Range("M99").Value = 605.2857
$ws.Range("H134").Value = 14480109
$ws.Range("I134").Value = 20001152
$ws.Range("J134").Value = 2977937
$ws.Range("K134").Value = 60003456
$ws.Range("L134").Value = 8933811
$ws.Range("M134").Value = -60000921
$ws.Range("N134").Value = -8938881

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 41674264
$ws.Range("J94").Value = 55565510
$ws.Range("L94").Value = 55565510
$ws.Range("N94").Value = -55566412
$ws.Range("H132").Value = 1707.3
$ws.Range("I132").Value = 1230.1892
$ws.Range("J132").Value = 3065.2307
$ws.Range("K132").Value = 3690.5676
$ws.Range("L132").Value = 9195.6921
$ws.Range("M132").Value = -1160.5676
$ws.Range("N132").Value = -14255.6921
$ws.Range("H134").Value = 873392.2
$ws.Range("I134").Value = 3864.7222
$ws.Range("J134").Value = 4003691
$ws.Range("K134").Value = 11594.1666
$ws.Range("L134").Value = 12011073
$ws.Range("M134").Value = -9059.1666
$ws.Range("N134").Value = -12016143

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 996.0417
$ws.Range("J98").Value = 1611.25
$ws.Range("L98").Value = 4833.75
$ws.Range("N98").Value = -7829.75
$ws.Range("H107").Value = 641280.4399999999
$ws.Range("I107").Value = 1508477.2
$ws.Range("J107").Value = 308.86957
$ws.Range("K107").Value = 4525431.6
$ws.Range("L107").Value = 926.60871
$ws.Range("M107").Value = -4523511.6
$ws.Range("N107").Value = -4766.60871

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 14706915
$ws.Range("I97").Value = 996.4167
$ws.Range("J97").Value = 50001120
$ws.Range("K97").Value = 996.4167
$ws.Range("L97").Value = 50001120
$ws.Range("M97").Value = -500.4167
$ws.Range("N97").Value = -50002112

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3000
$ws.Range("I46").Value = 3000
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -2812
$ws.Range("N46").ClearContents()
$ws.Range("H93").Value = 14398.85
$ws.Range("I93").Value = 3234.7856
$ws.Range("J93").Value = 40448.332
$ws.Range("K93").Value = 3234.7856
$ws.Range("L93").Value = 40448.332
$ws.Range("M93").Value = -1986.7856
$ws.Range("N93").Value = -42944.332
$ws.Range("H100").Value = 3116.2222
$ws.Range("I100").Value = 1222.7916
$ws.Range("K100").Value = 1222.7916
$ws.Range("M100").Value = -681.7916

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 25522.5
$ws.Range("J92").Value = 25522.5
$ws.Range("L92").Value = 25522.5
$ws.Range("N92").Value = -30514.5
$ws.Range("H93").Value = 29694.5
$ws.Range("J93").Value = 29694.5
$ws.Range("L93").Value = 29694.5
$ws.Range("N93").Value = -34686.5
$ws.Range("H94").Value = 14750
$ws.Range("J94").Value = 14750
$ws.Range("L94").Value = 14750
$ws.Range("N94").Value = -16552
$ws.Range("H95").Value = 19000
$ws.Range("J95").Value = 19000
$ws.Range("L95").Value = 19000
$ws.Range("N95").Value = -24492
$ws.Range("H96").Value = 1989.6842
$ws.Range("I96").Value = 1586.9333
$ws.Range("J96").Value = 3500
$ws.Range("K96").Value = 1586.9333
$ws.Range("L96").Value = 3500
$ws.Range("M96").Value = -213.9332999999999
$ws.Range("N96").Value = -6246
$ws.Range("H97").Value = 15000
$ws.Range("J97").Value = 15000
$ws.Range("L97").Value = 15000
$ws.Range("N97").Value = -16982
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H99").Value = 25158.666
$ws.Range("I99").Value = 20000
$ws.Range("J99").Value = 35476
$ws.Range("K99").Value = 20000
$ws.Range("L99").Value = 35476
$ws.Range("M99").Value = -17005
$ws.Range("N99").Value = -41466
$ws.Range("H100").Value = 2461.652
$ws.Range("I100").Value = 2564.4546
$ws.Range("J100").Value = 200
$ws.Range("K100").Value = 5128.9092
$ws.Range("L100").Value = 400
$ws.Range("M100").Value = -4587.9092
$ws.Range("N100").Value = -1482
$ws.Range("H101").Value = 13200.667
$ws.Range("J101").Value = 13200.667
$ws.Range("L101").Value = 13200.667
$ws.Range("N101").Value = -19690.667
$ws.Range("H102").Value = 35337
$ws.Range("J102").Value = 35337
$ws.Range("L102").Value = 35337
$ws.Range("N102").Value = -41827
$ws.Range("H103").Value = 31412.5
$ws.Range("J103").Value = 31412.5
$ws.Range("L103").Value = 31412.5
$ws.Range("N103").Value = -33756.5
$ws.Range("H104").Value = 50370
$ws.Range("J104").Value = 50370
$ws.Range("L104").Value = 50370
$ws.Range("N104").Value = -57358
$ws.Range("H105").Value = 30615
$ws.Range("J105").Value = 30615
$ws.Range("L105").Value = 30615
$ws.Range("N105").Value = -37603
$ws.Range("H106").Value = 35666
$ws.Range("J106").Value = 35666
$ws.Range("L106").Value = 35666
$ws.Range("N106").Value = -38190
$ws.Range("H136").Value = 7996.8438
$ws.Range("I136").Value = 6900.4736
$ws.Range("J136").Value = 9599.23
$ws.Range("K136").Value = 20701.4208
$ws.Range("L136").Value = 28797.69
$ws.Range("M136").Value = -18151.4208
$ws.Range("N136").Value = -33897.69

Write-Output "Applied Ridill_Profits scheduled data refresh."